# Generate Report for Handoff
# Adds two new handed-off files (29881203-... and 36f04959-...) as new rows
# to the "Overview", "zh-cn" and "de-de" worksheets, growing each table by
# two rows and refreshing the table/dimension ranges + hyperlinks.

$wb = $excel.ActiveWorkbook

$md1 = "29881203-4e87-4107-b19f-104f000f898d.md"
$md2 = "36f04959-9cc9-44ba-a8f1-ae4b14ccde34.md"
$path1 = "e2e\29881203-4e87-4107-b19f-104f000f898d.md"
$path2 = "e2e\36f04959-9cc9-44ba-a8f1-ae4b14ccde34.md"
$handoffDate = "2016-08-22 06:40:14"

$zhXlf1 = "29881203-4e87-4107-b19f-104f000f898d.fab48f9294dc927ab74bc1a4d08c23b9c4a76a3f.zh-cn.xlf"
$zhXlf2 = "36f04959-9cc9-44ba-a8f1-ae4b14ccde34.ad3e0f78246399f9759b5154fc3413ff53a00dfe.zh-cn.xlf"
$zhDate = "2016-08-22 06:40:00"

$deXlf1 = "29881203-4e87-4107-b19f-104f000f898d.fab48f9294dc927ab74bc1a4d08c23b9c4a76a3f.de-de.xlf"
$deXlf2 = "36f04959-9cc9-44ba-a8f1-ae4b14ccde34.ad3e0f78246399f9759b5154fc3413ff53a00dfe.de-de.xlf"

$status = "Ready for handoff"
$ext = ".md"
$srcPath = "e2e"
$priority = "ht"
$falseTxt = "False"
$trueTxt = "True"
$backDate = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview" -- two new rows (4 and 5), columns A:G
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Cells.Item(4, 1).Value = $md1
$wsOverview.Cells.Item(4, 2).Value = $path1
$wsOverview.Cells.Item(4, 3).Value = $ext
$wsOverview.Cells.Item(4, 4).Value = ""
$wsOverview.Cells.Item(4, 5).Value = $status
$wsOverview.Cells.Item(4, 6).Value = $status
$wsOverview.Cells.Item(4, 7).Value = $handoffDate

$wsOverview.Cells.Item(5, 1).Value = $md2
$wsOverview.Cells.Item(5, 2).Value = $path2
$wsOverview.Cells.Item(5, 3).Value = $ext
$wsOverview.Cells.Item(5, 4).Value = ""
$wsOverview.Cells.Item(5, 5).Value = $status
$wsOverview.Cells.Item(5, 6).Value = $status
$wsOverview.Cells.Item(5, 7).Value = $handoffDate

$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/29881203-4e87-4107-b19f-104f000f898d.md", "", "", $path1) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/36f04959-9cc9-44ba-a8f1-ae4b14ccde34.md", "", "", $path2) | Out-Null

$loOverview = $wsOverview.ListObjects.Item("Overview")
$loOverview.Resize($wsOverview.Range("A1:G5"))

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn" -- two new rows (4 and 5), columns A:P
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Cells.Item(4, 1).Value = $md1
$wsZh.Cells.Item(4, 2).Value = $ext
$wsZh.Cells.Item(4, 3).Value = $status
$wsZh.Cells.Item(4, 4).Value = $srcPath
$wsZh.Cells.Item(4, 5).Value = $priority
$wsZh.Cells.Item(4, 6).Value = $falseTxt
$wsZh.Cells.Item(4, 7).Value = $zhXlf1
$wsZh.Cells.Item(4, 8).Value = $zhDate
$wsZh.Cells.Item(4, 9).Value = ""
$wsZh.Cells.Item(4, 10).Value = ""
$wsZh.Cells.Item(4, 11).Value = $backDate
$wsZh.Cells.Item(4, 12).Value = ""
$wsZh.Cells.Item(4, 13).Value = $trueTxt
$wsZh.Cells.Item(4, 14).Value = ""
$wsZh.Cells.Item(4, 15).Value = $falseTxt
$wsZh.Cells.Item(4, 16).Value = ""

$wsZh.Cells.Item(5, 1).Value = $md2
$wsZh.Cells.Item(5, 2).Value = $ext
$wsZh.Cells.Item(5, 3).Value = $status
$wsZh.Cells.Item(5, 4).Value = $srcPath
$wsZh.Cells.Item(5, 5).Value = $priority
$wsZh.Cells.Item(5, 6).Value = $falseTxt
$wsZh.Cells.Item(5, 7).Value = $zhXlf2
$wsZh.Cells.Item(5, 8).Value = $zhDate
$wsZh.Cells.Item(5, 9).Value = ""
$wsZh.Cells.Item(5, 10).Value = ""
$wsZh.Cells.Item(5, 11).Value = $backDate
$wsZh.Cells.Item(5, 12).Value = ""
$wsZh.Cells.Item(5, 13).Value = $trueTxt
$wsZh.Cells.Item(5, 14).Value = ""
$wsZh.Cells.Item(5, 15).Value = $falseTxt
$wsZh.Cells.Item(5, 16).Value = ""

$wsZh.Hyperlinks.Add($wsZh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/29881203-4e87-4107-b19f-104f000f898d.md", "", "", $md1) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/36f04959-9cc9-44ba-a8f1-ae4b14ccde34.md", "", "", $md2) | Out-Null

$loZh = $wsZh.ListObjects.Item("zh-cn")
$loZh.Resize($wsZh.Range("A1:P5"))

# ---------------------------------------------------------------------------
# Sheet 3: "de-de" -- two new rows (4 and 5), columns A:P
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Cells.Item(4, 1).Value = $md1
$wsDe.Cells.Item(4, 2).Value = $ext
$wsDe.Cells.Item(4, 3).Value = $status
$wsDe.Cells.Item(4, 4).Value = $srcPath
$wsDe.Cells.Item(4, 5).Value = $priority
$wsDe.Cells.Item(4, 6).Value = $falseTxt
$wsDe.Cells.Item(4, 7).Value = $deXlf1
$wsDe.Cells.Item(4, 8).Value = $handoffDate
$wsDe.Cells.Item(4, 9).Value = ""
$wsDe.Cells.Item(4, 10).Value = ""
$wsDe.Cells.Item(4, 11).Value = $backDate
$wsDe.Cells.Item(4, 12).Value = ""
$wsDe.Cells.Item(4, 13).Value = $trueTxt
$wsDe.Cells.Item(4, 14).Value = ""
$wsDe.Cells.Item(4, 15).Value = $falseTxt
$wsDe.Cells.Item(4, 16).Value = ""

$wsDe.Cells.Item(5, 1).Value = $md2
$wsDe.Cells.Item(5, 2).Value = $ext
$wsDe.Cells.Item(5, 3).Value = $status
$wsDe.Cells.Item(5, 4).Value = $srcPath
$wsDe.Cells.Item(5, 5).Value = $priority
$wsDe.Cells.Item(5, 6).Value = $falseTxt
$wsDe.Cells.Item(5, 7).Value = $deXlf2
$wsDe.Cells.Item(5, 8).Value = $handoffDate
$wsDe.Cells.Item(5, 9).Value = ""
$wsDe.Cells.Item(5, 10).Value = ""
$wsDe.Cells.Item(5, 11).Value = $backDate
$wsDe.Cells.Item(5, 12).Value = ""
$wsDe.Cells.Item(5, 13).Value = $trueTxt
$wsDe.Cells.Item(5, 14).Value = ""
$wsDe.Cells.Item(5, 15).Value = $falseTxt
$wsDe.Cells.Item(5, 16).Value = ""

$wsDe.Hyperlinks.Add($wsDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/29881203-4e87-4107-b19f-104f000f898d.md", "", "", $md1) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/master/e2e/36f04959-9cc9-44ba-a8f1-ae4b14ccde34.md", "", "", $md2) | Out-Null

$loDe = $wsDe.ListObjects.Item("de-de")
$loDe.Resize($wsDe.Range("A1:P5"))

Write-Host "Report generated for handback: added 2 rows to Overview, zh-cn, de-de"
